$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Row 2 (Beta) ---
$ws.Range("C2").Value = 40.352730803732612230305676348507
$ws.Range("D2").Value = 0.000743282006413391576477023293
$ws.Range("E2").Value = 0.057689018187517870916991569175
$ws.Range("F2").Value = 40.664306090179252350935712456703
$ws.Range("G2").Value = 39.879411435214272785287903388962
$ws.Range("H2").Value = 41.441740667327252367613255046308
$ws.Range("I2").Value = 0.001004674904283017929096377685
$ws.Range("J2").Value = 0.000742514493211707043637814785
$ws.Range("K2").Value = 0.001505250748018314033843667588
$ws.Range("L2").Value = 0.059168268540023619128653820098
$ws.Range("M2").Value = 0.058128481506161337744043038356
$ws.Range("N2").Value = 0.060220047660702338832017233017

# --- Update Row 3 (Gamma) ---
$ws.Range("F3").Value = 0.000009923077592035762519405429
$ws.Range("G3").Value = 0.000000002056593222115246101122
$ws.Range("H3").Value = 0.000028981311265019180561713572
$ws.Range("I3").Value = 0.000008334097948465038171225044
$ws.Range("J3").Value = 0.000000001915140631978846994847
$ws.Range("K3").Value = 0.000024277609599195600298621525
$ws.Range("L3").Value = 0.000010177236962385319290645344
$ws.Range("M3").Value = 0.000000002136343402779747826248
$ws.Range("N3").Value = 0.000029716767113971071022186787

# --- Add Row 4 (Beta + Gamma) ---
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 40.352730803732612230305676348507
$ws.Range("D4").Value = 0.000743282006413391576477023293
$ws.Range("E4").Value = 0.057689018187517870916991569175
$ws.Range("F4").Value = 40.664316013256836868094978854060
$ws.Range("G4").Value = 39.879411437270860574244579765946
$ws.Range("H4").Value = 41.441769648638519640826416434720
$ws.Range("I4").Value = 0.001013009002231483055359029244
$ws.Range("J4").Value = 0.000742516408352338953303972691
$ws.Range("K4").Value = 0.001529528357617510091540080630
$ws.Range("L4").Value = 0.059178445776985999504660185266
$ws.Range("M4").Value = 0.058128483642504752415547386590
$ws.Range("N4").Value = 0.060249764427816307477137058868
